$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "Golang Architect / Principal Backend Architect"
$ws.Range("B72").Value = "https://www.dice.com/job-detail/82042fb2-5faf-481f-b812-13a7e928975b"
$ws.Range("C72").Value = "Atlanta, Georgia"
$ws.Range("D72").Value = "Third Party"
$ws.Range("E72").Value = "$65 - $75 per hour"
$ws.Range("F72").Value = "STAND 8"
